$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2966.889
$ws.Range("I74").Value = 2675.75
$ws.Range("K74").Value = 2675.75
$ws.Range("M74").Value = -1739.75
$ws.Range("H76").Value = 3623.0435
$ws.Range("I76").Value = 3746.111
$ws.Range("J76").Value = 3180
$ws.Range("K76").Value = 3746.111
$ws.Range("L76").Value = 3180
$ws.Range("M76").Value = -3431.111
$ws.Range("N76").Value = -3810
$ws.Range("H77").Value = 2966.889
$ws.Range("I77").Value = 2675.75
$ws.Range("K77").Value = 13378.75
$ws.Range("M77").Value = -8698.75
$ws.Range("H79").Value = 3623.0435
$ws.Range("I79").Value = 3746.111
$ws.Range("J79").Value = 3180
$ws.Range("K79").Value = 3746.111
$ws.Range("L79").Value = 3180
$ws.Range("M79").Value = -2654.111
$ws.Range("N79").Value = -5364
$ws.Range("H111").Value = 2869.5386
$ws.Range("I111").Value = 2824.8572
$ws.Range("J111").Value = 2921.6667
$ws.Range("K111").Value = 8474.571599999999
$ws.Range("L111").Value = 8765.000100000001
$ws.Range("M111").Value = -5407.571599999999
$ws.Range("N111").Value = -14899.0001
$ws.Range("H132").Value = 7252542
$ws.Range("I132").Value = 11117218
$ws.Range("J132").Value = 6273.9375
$ws.Range("K132").Value = 33351654
$ws.Range("L132").Value = 18821.8125
$ws.Range("M132").Value = -33349124
$ws.Range("N132").Value = -23881.8125
$ws.Range("H137").Value = 1237.9354
$ws.Range("I137").Value = 906.4857
$ws.Range("J137").Value = 1667.5927
$ws.Range("K137").Value = 2719.4571
$ws.Range("L137").Value = 5002.7781
$ws.Range("M137").Value = -169.4570999999996
$ws.Range("N137").Value = -10102.7781
$ws.Range("H138").Value = 525319.25
$ws.Range("J138").Value = 837883.1
$ws.Range("L138").Value = 2513649.3
$ws.Range("N138").Value = -2523929.3
$ws.Range("H141").Value = 1012.8571
$ws.Range("I141").Value = 848.3333
$ws.Range("K141").Value = 2544.9999
$ws.Range("M141").Value = 2635.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3970.6128
$ws.Range("I32").Value = 3578.5066
$ws.Range("J32").Value = 5857.625
$ws.Range("K32").Value = 3578.5066
$ws.Range("L32").Value = 5857.625
$ws.Range("M32").Value = -3291.5066
$ws.Range("N32").Value = -6431.625
$ws.Range("H61").Value = 47620468
$ws.Range("I61").Value = 62500930
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 62500930
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -62500718
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 1228.4286
$ws.Range("I74").Value = 976.2308
$ws.Range("K74").Value = 976.2308
$ws.Range("M74").Value = -102.2308
$ws.Range("H77").Value = 1228.4286
$ws.Range("I77").Value = 976.2308
$ws.Range("K77").Value = 4881.154
$ws.Range("M77").Value = -513.1540000000005
$ws.Range("H132").Value = 1868.25
$ws.Range("I132").Value = 1561.4348
$ws.Range("J132").Value = 3279.6
$ws.Range("K132").Value = 4684.3044
$ws.Range("L132").Value = 9838.799999999999
$ws.Range("M132").Value = -2154.3044
$ws.Range("N132").Value = -14898.8
$ws.Range("H136").Value = 47620468
$ws.Range("I136").Value = 62500930
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 187502790
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -187500240
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3373.0425
$ws.Range("I134").Value = 1047.9736
$ws.Range("J134").Value = 13190
$ws.Range("K134").Value = 3143.9208
$ws.Range("L134").Value = 39570
$ws.Range("M134").Value = -608.9207999999999
$ws.Range("N134").Value = -44640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2465.8096
$ws.Range("I31").Value = 2611.375
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 2611.375
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -2316.375
$ws.Range("N31").Value = -2590
$ws.Range("H34").Value = 2465.8096
$ws.Range("I34").Value = 2611.375
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 2611.375
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -2409.375
$ws.Range("N34").Value = -2404
$ws.Range("H51").Value = 24000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 24000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 24000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -25472
$ws.Range("H61").Value = 24000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 24000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 24000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -24696
$ws.Range("H99").Value = 1254675.5
$ws.Range("I99").Value = 1755785
$ws.Range("J99").Value = 1901.6666
$ws.Range("K99").Value = 1755785
$ws.Range("L99").Value = 1901.6666
$ws.Range("M99").Value = -1754287
$ws.Range("N99").Value = -4897.6666
$ws.Range("H126").Value = 1254675.5
$ws.Range("I126").Value = 1755785
$ws.Range("J126").Value = 1901.6666
$ws.Range("K126").Value = 5267355
$ws.Range("L126").Value = 5704.9998
$ws.Range("M126").Value = -5264885
$ws.Range("N126").Value = -10644.9998
$ws.Range("H132").Value = 2252.1428
$ws.Range("I132").Value = 1962.5883
$ws.Range("K132").Value = 5887.7649
$ws.Range("M132").Value = -3357.7649
$ws.Range("I134").Value = 1212.9333
$ws.Range("J134").Value = 55556900
$ws.Range("K134").Value = 3638.7999
$ws.Range("L134").Value = 166670700
$ws.Range("M134").Value = -1103.7999
$ws.Range("N134").Value = -166675770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 682.78125
$ws.Range("I113").Value = 499.2
$ws.Range("J113").Value = 716.7778
$ws.Range("K113").Value = 1497.6
$ws.Range("L113").Value = 2150.3334
$ws.Range("M113").Value = 672.4000000000001
$ws.Range("N113").Value = -6490.3334
$ws.Range("H131").Value = 38466520
$ws.Range("I131").Value = 200000480
$ws.Range("J131").Value = 6052.3335
$ws.Range("K131").Value = 600001440
$ws.Range("L131").Value = 18157.0005
$ws.Range("M131").Value = -599996400
$ws.Range("N131").Value = -28237.0005
$ws.Range("H132").Value = 1433.3334
$ws.Range("I132").Value = 1300
$ws.Range("K132").Value = 11700
$ws.Range("M132").Value = -9170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 24900
$ws.Range("J26").Value = 24900
$ws.Range("L26").Value = 24900
$ws.Range("N26").Value = -25460
$ws.Range("H50").Value = 24900
$ws.Range("J50").Value = 24900
$ws.Range("L50").Value = 24900
$ws.Range("N50").Value = -25896
$ws.Range("H54").Value = 3090
$ws.Range("J54").Value = 3090
$ws.Range("L54").Value = 3090
$ws.Range("N54").Value = -3870
$ws.Range("H132").Value = 3512.3
$ws.Range("I132").Value = 3111.3845
$ws.Range("K132").Value = 9334.1535
$ws.Range("M132").Value = -6804.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1607.5
$ws.Range("I7").Value = 1607.5
$ws.Range("K7").Value = 1607.5
$ws.Range("M7").Value = -1495.5
$ws.Range("H40").Value = 2491.16
$ws.Range("I40").Value = 1830.4286
$ws.Range("K40").Value = 1830.4286
$ws.Range("M40").Value = -1694.4286
$ws.Range("H122").Value = 14170594
$ws.Range("I122").Value = 28335022
$ws.Range("K122").Value = 85005066
$ws.Range("M122").Value = -85002616
$ws.Range("H126").Value = 1607.5
$ws.Range("I126").Value = 1607.5
$ws.Range("K126").Value = 4822.5
$ws.Range("M126").Value = -2352.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 475.07144
$ws.Range("I113").Value = 426.2
$ws.Range("J113").Value = 502.22223
$ws.Range("K113").Value = 1278.6
$ws.Range("L113").Value = 1506.66669
$ws.Range("M113").Value = 891.4000000000001
$ws.Range("N113").Value = -5846.66669
$ws.Range("H126").Value = 111112856
$ws.Range("I126").Value = 123458390
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 370375170
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -370372700
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 2742
$ws.Range("I132").Value = 2799.9
$ws.Range("K132").Value = 8399.700000000001
$ws.Range("M132").Value = -5869.700000000001
$ws.Range("H136").Value = 721.7143
$ws.Range("I136").Value = 646.3333
$ws.Range("J136").Value = 947.8570999999999
$ws.Range("K136").Value = 1938.9999
$ws.Range("L136").Value = 2843.5713
$ws.Range("M136").Value = 611.0001
$ws.Range("N136").Value = -7943.5713
